$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I20").Value = -1.417434297872755
$ws.Range("J20").Value = 0.2306493530254091
$ws.Range("K20").Value = 0.2501630079354575
$ws.Range("L20").Value = 2.346422858371362
